$d = $word.ActiveDocument

# Pull the full package OOXML so we can do precise, low-level surgery on
# attributes (wp:docPr/@id, wp:docPr/@name, pic:cNvPr/@id, pic:cNvPr/@name)
# that are not exposed as settable properties on InlineShape in the Word
# object model.
$range = $d.Content
$xml = $range.WordOpenXML

# --- Step 1: retarget the wp:docPr id/name pairs (these strings are each
# unique in the package, so plain substring replacement is safe/unambiguous).
$docPrMap = @{
    '<wp:docPr id="34" name="Picture 34"/>' = '<wp:docPr id="18" name="Picture 18"/>'
    '<wp:docPr id="35" name="Picture 35"/>' = '<wp:docPr id="19" name="Picture 19"/>'
    '<wp:docPr id="36" name="Picture 36"/>' = '<wp:docPr id="20" name="Picture 20"/>'
    '<wp:docPr id="37" name="Picture 37"/>' = '<wp:docPr id="21" name="Picture 21"/>'
    '<wp:docPr id="38" name="Picture 38"/>' = '<wp:docPr id="22" name="Picture 22"/>'
    '<wp:docPr id="39" name="Picture 39"/>' = '<wp:docPr id="23" name="Picture 23"/>'
    '<wp:docPr id="40" name="Picture 40"/>' = '<wp:docPr id="24" name="Picture 24"/>'
    '<wp:docPr id="41" name="Picture 41"/>' = '<wp:docPr id="25" name="Picture 25"/>'
    '<wp:docPr id="42" name="Picture 42"/>' = '<wp:docPr id="26" name="Picture 26"/>'
    '<wp:docPr id="43" name="Picture 43"/>' = '<wp:docPr id="27" name="Picture 27"/>'
    '<wp:docPr id="44" name="Picture 44"/>' = '<wp:docPr id="28" name="Picture 28"/>'
    '<wp:docPr id="45" name="Picture 45"/>' = '<wp:docPr id="29" name="Picture 29"/>'
    '<wp:docPr id="46" name="Picture 46"/>' = '<wp:docPr id="30" name="Picture 30"/>'
    '<wp:docPr id="47" name="Picture 47"/>' = '<wp:docPr id="31" name="Picture 31"/>'
    '<wp:docPr id="48" name="Picture 48"/>' = '<wp:docPr id="32" name="Picture 32"/>'
    '<wp:docPr id="49" name="Picture 49"/>' = '<wp:docPr id="33" name="Picture 33"/>'
}

foreach ($old in $docPrMap.Keys) {
    $xml = $xml.Replace($old, $docPrMap[$old])
}

# --- Step 2: retarget the pic:cNvPr name for each picture (id stays "0" in
# every case). All 16 occurrences currently share the literal text
# '<pic:cNvPr id="0" name="image.png"/>', so they must be substituted in
# the order they occur in the document (which matches the wp:docPr order
# above: 18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33).
$newNames = @(
    "Leo-0_05.png",
    "Leo-0_15.png",
    "Leo-0_25.png",
    "Leo-0_35.png",
    "Leo-0_45.png",
    "Leo-0_55.png",
    "Leo-0_65.png",
    "Leo-0_75.png",
    "Leo-0_05.png",
    "Leo-0_15.png",
    "Leo-0_25.png",
    "Leo-0_35.png",
    "Leo-0_45.png",
    "Leo-0_55.png",
    "Leo-0_65.png",
    "Leo-0_75.png"
)

$search = '<pic:cNvPr id="0" name="image.png"/>'
$rebuilt = ""
$remaining = $xml
$i = 0
while ($true) {
    $pos = $remaining.IndexOf($search)
    if ($pos -lt 0 -or $i -ge $newNames.Length) {
        $rebuilt += $remaining
        break
    }
    $rebuilt += $remaining.Substring(0, $pos)
    $rebuilt += '<pic:cNvPr id="0" name="' + $newNames[$i] + '"/>'
    $remaining = $remaining.Substring($pos + $search.Length)
    $i++
}
$xml = $rebuilt

# --- Step 3: write the package back.
$range.WordOpenXML = $xml

Write-Output "replaced docPr entries and $i cNvPr names"
